$wb = $excel.ActiveWorkbook

# --- "Alle Abteilungen": literal respondent count for "4. Semester" /
#     Landschaftsarchitektur (H4) goes from 1 to 6 -----------------------
$wsAlle = $wb.Worksheets.Item("Alle Abteilungen")
$wsAlle.Range("H4").Value = 6

# --- "Landschaftsarchitektur": five newly-entered questionnaires for the
#     "4. Semester" block (rows 26-33, columns C-G) ----------------------
$wsLA = $wb.Worksheets.Item("Landschaftsarchitektur")

$wsLA.Range("C26").Value = 3
$wsLA.Range("D26").Value = 4
$wsLA.Range("E26").Value = 4
$wsLA.Range("F26").Value = 3
$wsLA.Range("G26").Value = 3

$wsLA.Range("C27").Value = 2
$wsLA.Range("D27").Value = 2
$wsLA.Range("E27").Value = 4
$wsLA.Range("F27").Value = 2
$wsLA.Range("G27").Value = 1

$wsLA.Range("C28").Value = 4
$wsLA.Range("D28").Value = 0
$wsLA.Range("E28").Value = 4
$wsLA.Range("F28").Value = 0
$wsLA.Range("G28").Value = 3

$wsLA.Range("C29").Value = 2
$wsLA.Range("D29").Value = 3
$wsLA.Range("E29").Value = 3
$wsLA.Range("F29").Value = 4
$wsLA.Range("G29").Value = 1

$wsLA.Range("C30").Value = 3
$wsLA.Range("D30").Value = 4
$wsLA.Range("E30").Value = 3
$wsLA.Range("F30").Value = 4
$wsLA.Range("G30").Value = 4

$wsLA.Range("C31").Value = 3
$wsLA.Range("D31").Value = 0
$wsLA.Range("E31").Value = 3
$wsLA.Range("F31").Value = 2
$wsLA.Range("G31").Value = 1

$wsLA.Range("C32").Value = 2
$wsLA.Range("D32").Value = 2
$wsLA.Range("E32").Value = 2
$wsLA.Range("F32").Value = 1
$wsLA.Range("G32").Value = 1

$wsLA.Range("C33").Value = 3
$wsLA.Range("D33").Value = 2
$wsLA.Range("E33").Value = 3
$wsLA.Range("F33").Value = 1
$wsLA.Range("G33").Value = 2

# --- Recalculate so every dependent formula (and all the embedded charts
#     that read from these ranges) reflects the new figures --------------
$excel.CalculateFullRebuild()

# --- Restore the selections recorded in each sheet's view, activating
#     "Alle Abteilungen" last so it remains the selected tab -------------
$wsLA.Activate()
$wsLA.Range("G34").Select()

$wsAlle.Activate()
$wsAlle.Range("H5").Select()
